$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.317308670475938
$ws.Cells.Item(2, 3).Value = 0.3856787182343169
$ws.Cells.Item(2, 4).Value = 0.2869411618068085
$ws.Cells.Item(2, 6).Value = 1.41812552750784
$ws.Cells.Item(2, 7).Value = 0.7526355808948892
$ws.Cells.Item(2, 8).Value = 0.8142458949535722
$ws.Cells.Item(2, 9).Value = 0.5428215978979907
$ws.Cells.Item(2, 10).Value = 0.3487309208576619

$ws.Cells.Item(3, 2).Value = 1.176388625456809
$ws.Cells.Item(3, 3).Value = 0.3405951396038631
$ws.Cells.Item(3, 4).Value = 0.2786549291283222
$ws.Cells.Item(3, 6).Value = 1.41306692097649
$ws.Cells.Item(3, 7).Value = 0.749790072561467
$ws.Cells.Item(3, 8).Value = 0.8197619046344897
$ws.Cells.Item(3, 9).Value = 0.5526300047420065
$ws.Cells.Item(3, 10).Value = 0.3373669850378178

$ws.Cells.Item(4, 2).Value = 1.089869530970589
$ws.Cells.Item(4, 3).Value = 0.3129009597249137
$ws.Cells.Item(4, 4).Value = 0.2736526380471958
$ws.Cells.Item(4, 6).Value = 1.4111786068149
$ws.Cells.Item(4, 7).Value = 0.7488929435580474
$ws.Cells.Item(4, 8).Value = 0.8237784926231626
$ws.Cells.Item(4, 9).Value = 0.559213172032063
$ws.Cells.Item(4, 10).Value = 0.3306071877485977

$ws.Cells.Item(5, 2).Value = 1.054615097239548
$ws.Cells.Item(5, 3).Value = 0.3016123578272811
$ws.Cells.Item(5, 4).Value = 0.2716357505789659
$ws.Cells.Item(5, 6).Value = 1.410714054655571
$ws.Cells.Item(5, 7).Value = 0.7487397978817967
$ws.Cells.Item(5, 8).Value = 0.825573166420952
$ws.Cells.Item(5, 9).Value = 0.5620362697893491
$ws.Cells.Item(5, 10).Value = 0.3279070264907631

$ws.Cells.Item(6, 2).Value = 1.048761330440755
$ws.Cells.Item(6, 3).Value = 0.2997377188143844
$ws.Cells.Item(6, 4).Value = 0.271302154433954
$ws.Cells.Item(6, 6).Value = 1.410655298007413
$ws.Cells.Item(6, 7).Value = 0.7487271609718391
$ws.Cells.Item(6, 8).Value = 0.8258806925484095
$ws.Cells.Item(6, 9).Value = 0.562513505331502
$ws.Cells.Item(6, 10).Value = 0.3274619523644589

$ws.Cells.Item(7, 2).Value = 1.089394064402256
$ws.Cells.Item(7, 3).Value = 0.312748729571382
$ws.Cells.Item(7, 4).Value = 0.2736253500724786
$ws.Cells.Item(7, 6).Value = 1.411171108626306
$ws.Cells.Item(7, 7).Value = 0.7488900197462556
$ws.Cells.Item(7, 8).Value = 0.8238020575964384
$ws.Cells.Item(7, 9).Value = 0.5592506776348429
$ws.Cells.Item(7, 10).Value = 0.3305705520199922

$ws.Cells.Item(8, 2).Value = 1.268718776771721
$ws.Cells.Item(8, 3).Value = 0.3701366290666215
$ws.Cells.Item(8, 4).Value = 0.2840663756371242
$ws.Cells.Item(8, 6).Value = 1.416127830767266
$ws.Cells.Item(8, 7).Value = 0.7514772386483912
$ws.Cells.Item(8, 8).Value = 0.8160168288635816
$ws.Cells.Item(8, 9).Value = 0.5460867860007106
$ws.Cells.Item(8, 10).Value = 0.344767305971331

$ws.Cells.Item(9, 2).Value = 1.620391380039734
$ws.Cells.Item(9, 3).Value = 0.4825713072108329
$ws.Cells.Item(9, 4).Value = 0.3052169180657671
$ws.Cells.Item(9, 6).Value = 1.435568939548858
$ws.Cells.Item(9, 7).Value = 0.763355704855968
$ws.Cells.Item(9, 8).Value = 0.8057685201907958
$ws.Cells.Item(9, 9).Value = 0.5247489195863793
$ws.Cells.Item(9, 10).Value = 0.3743469224611857

$ws.Cells.Item(10, 2).Value = 1.878759560914204
$ws.Cells.Item(10, 3).Value = 0.5651205014091829
$ws.Cells.Item(10, 4).Value = 0.3211666670730153
$ws.Cells.Item(10, 6).Value = 1.455864617838074
$ws.Cells.Item(10, 7).Value = 0.776317116068725
$ws.Cells.Item(10, 8).Value = 0.8013296025412444
$ws.Cells.Item(10, 9).Value = 0.5118380751502194
$ws.Cells.Item(10, 10).Value = 0.3971597152056177

$ws.Cells.Item(11, 2).Value = 1.996295612751624
$ws.Cells.Item(11, 3).Value = 0.6026640325525818
$ws.Cells.Item(11, 4).Value = 0.3285115954543585
$ws.Cells.Item(11, 6).Value = 1.466422050950555
$ws.Cells.Item(11, 7).Value = 0.7831516919508914
$ws.Cells.Item(11, 8).Value = 0.7999880981716814
$ws.Cells.Item(11, 9).Value = 0.5065732239845531
$ws.Cells.Item(11, 10).Value = 0.407777023168677

$ws.Cells.Item(12, 2).Value = 2.040803257151708
$ws.Cells.Item(12, 3).Value = 0.6168795592290621
$ws.Cells.Item(12, 4).Value = 0.331305721010807
$ws.Cells.Item(12, 6).Value = 1.470611848311563
$ws.Cells.Item(12, 7).Value = 0.7858761784986115
$ws.Cells.Item(12, 8).Value = 0.7995781154668862
$ws.Cells.Item(12, 9).Value = 0.5046677094756902
$ws.Cells.Item(12, 10).Value = 0.4118322980615261

$ws.Cells.Item(13, 2).Value = 2.031217788941547
$ws.Cells.Item(13, 3).Value = 0.6138180557810529
$ws.Cells.Item(13, 4).Value = 0.3307033898225313
$ws.Cells.Item(13, 6).Value = 1.469700941837004
$ws.Cells.Item(13, 7).Value = 0.7852833229261478
$ws.Cells.Item(13, 8).Value = 0.7996620444742888
$ws.Cells.Item(13, 9).Value = 0.5050741637144398
$ws.Cells.Item(13, 10).Value = 0.4109573730746092

$ws.Cells.Item(14, 2).Value = 1.999957304068573
$ws.Cells.Item(14, 3).Value = 0.6038335813522622
$ws.Cells.Item(14, 4).Value = 0.3287412145633652
$ws.Cells.Item(14, 6).Value = 1.466762893458409
$ws.Cells.Item(14, 7).Value = 0.783373095525647
$ws.Cells.Item(14, 8).Value = 0.7999524010926251
$ws.Cells.Item(14, 9).Value = 0.5064146851368463
$ws.Cells.Item(14, 10).Value = 0.408109955464937

$ws.Cells.Item(15, 2).Value = 1.980809228862768
$ws.Cells.Item(15, 3).Value = 0.5977176161180182
$ws.Cells.Item(15, 4).Value = 0.3275409849842958
$ws.Cells.Item(15, 6).Value = 1.464988290707254
$ws.Cells.Item(15, 7).Value = 0.7822208298907753
$ws.Cells.Item(15, 8).Value = 0.800143033893761
$ws.Cells.Item(15, 9).Value = 0.5072472960983987
$ws.Cells.Item(15, 10).Value = 0.406370360978201

$ws.Cells.Item(16, 2).Value = 1.871078266188192
$ws.Cells.Item(16, 3).Value = 0.5626667572996098
$ws.Cells.Item(16, 4).Value = 0.3206884487804302
$ws.Cells.Item(16, 6).Value = 1.455201444715442
$ws.Cells.Item(16, 7).Value = 0.7758894662005531
$ws.Cells.Item(16, 8).Value = 0.8014309587960895
$ws.Cells.Item(16, 9).Value = 0.5121944461287029
$ws.Cells.Item(16, 10).Value = 0.39647069450983

$ws.Cells.Item(17, 2).Value = 1.803761831422321
$ws.Cells.Item(17, 3).Value = 0.5411618862318051
$ws.Cells.Item(17, 4).Value = 0.3165074545952393
$ws.Cells.Item(17, 6).Value = 1.449537858359562
$ws.Cells.Item(17, 7).Value = 0.7722467442715129
$ws.Cells.Item(17, 8).Value = 0.8023950466957217
$ws.Cells.Item(17, 9).Value = 0.5153856669865426
$ws.Cells.Item(17, 10).Value = 0.3904591410582867

$ws.Cells.Item(18, 2).Value = 1.765043553848784
$ws.Cells.Item(18, 3).Value = 0.5287920614629229
$ws.Cells.Item(18, 4).Value = 0.3141110680948032
$ws.Cells.Item(18, 6).Value = 1.446404927118721
$ws.Cells.Item(18, 7).Value = 0.7702397836610402
$ws.Cells.Item(18, 8).Value = 0.8030133274124012
$ws.Cells.Item(18, 9).Value = 0.5172783929020568
$ws.Cells.Item(18, 10).Value = 0.3870239955983408

$ws.Cells.Item(19, 2).Value = 1.75193430699693
$ws.Cells.Item(19, 3).Value = 0.5246037207302834
$ws.Cells.Item(19, 4).Value = 0.31330113989344
$ws.Cells.Item(19, 6).Value = 1.445365529876412
$ws.Cells.Item(19, 7).Value = 0.7695753720495588
$ws.Cells.Item(19, 8).Value = 0.8032336012023507
$ws.Cells.Item(19, 9).Value = 0.5179290440029654
$ws.Cells.Item(19, 10).Value = 0.3858647776923476

$ws.Cells.Item(20, 2).Value = 1.810927752488681
$ws.Cells.Item(20, 3).Value = 0.5434511996643323
$ws.Cells.Item(20, 4).Value = 0.3169516586792724
$ws.Cells.Item(20, 6).Value = 1.450127849519305
$ws.Cells.Item(20, 7).Value = 0.7726253756418657
$ws.Cells.Item(20, 8).Value = 0.802285814713926
$ws.Cells.Item(20, 9).Value = 0.5150400290480839
$ws.Cells.Item(20, 10).Value = 0.3910967463298931

$ws.Cells.Item(21, 2).Value = 2.009139294411739
$ws.Cells.Item(21, 3).Value = 0.6067663026418586
$ws.Cells.Item(21, 4).Value = 0.3293172073207131
$ws.Cells.Item(21, 6).Value = 1.467620649316075
$ws.Cells.Item(21, 7).Value = 0.7839304627349009
$ws.Cells.Item(21, 8).Value = 0.7998644517791575
$ws.Cells.Item(21, 9).Value = 0.506018543075875
$ws.Cells.Item(21, 10).Value = 0.408945366294688

$ws.Cells.Item(22, 2).Value = 2.138677365033004
$ws.Cells.Item(22, 3).Value = 0.6481381880795425
$ws.Cells.Item(22, 4).Value = 0.3374731505768409
$ws.Cells.Item(22, 6).Value = 1.48017256358969
$ws.Cells.Item(22, 7).Value = 0.7921145976313824
$ws.Cells.Item(22, 8).Value = 0.7988535210265582
$ws.Cells.Item(22, 9).Value = 0.5006367264603497
$ws.Cells.Item(22, 10).Value = 0.4208130363119267

$ws.Cells.Item(23, 2).Value = 2.069541252210399
$ws.Cells.Item(23, 3).Value = 0.6260580370827142
$ws.Cells.Item(23, 4).Value = 0.3331133941530595
$ws.Cells.Item(23, 6).Value = 1.47337047887612
$ws.Cells.Item(23, 7).Value = 0.7876732939475488
$ws.Cells.Item(23, 8).Value = 0.7993405920069989
$ws.Cells.Item(23, 9).Value = 0.5034618213455886
$ws.Cells.Item(23, 10).Value = 0.4144604129516409

$ws.Cells.Item(24, 2).Value = 1.807688092948638
$ws.Cells.Item(24, 3).Value = 0.5424162209509973
$ws.Cells.Item(24, 4).Value = 0.3167508111981476
$ws.Cells.Item(24, 6).Value = 1.449860731054699
$ws.Cells.Item(24, 7).Value = 0.7724539246146662
$ws.Cells.Item(24, 8).Value = 0.8023349991356667
$ws.Cells.Item(24, 9).Value = 0.5151961111891339
$ws.Cells.Item(24, 10).Value = 0.3908084196709467

$ws.Cells.Item(25, 2).Value = 1.525254590027316
$ws.Cells.Item(25, 3).Value = 0.4521651661479495
$ws.Cells.Item(25, 4).Value = 0.299422959233226
$ws.Cells.Item(25, 6).Value = 1.429260027421961
$ws.Cells.Item(25, 7).Value = 0.759404892005449
$ws.Cells.Item(25, 8).Value = 0.8080004893578092
$ws.Cells.Item(25, 9).Value = 0.5300383116024712
$ws.Cells.Item(25, 10).Value = 0.3661565370953213

Write-Host "Updated pl_mw values for case with 380 kV"
